$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of 2021-season results, appended below the existing data in the
# same layout already used on the sheet: Year, Manager, Finish, Wins,
# Points, PointsAllowed, Playoffs (columns A-G). Columns H/I (DraftPosition /
# Bye) are intentionally left blank, matching the source rows, and column J
# keeps the "Champion" helper formula going.
$newRows = @(
    @(2021, "Colin",    3,    8, 1894.1,              1747.52, 1),
    @(2021, "John",     4,    7, 1926.56,             1853.62, 1),
    @(2021, "Charles",  9,    6, 1669.2,              1772.42, 0),
    @(2021, "Jennifer", 8,    7, 1733.16,             1712.72, 0),
    @(2021, "Chester",  10,   4, 1660.82,             1734.22, 0),
    @(2021, "EricR",    7,    7, 1719.06,             1691.04, 0),
    @(2021, "Mike",     6,    8, 1736.96,             1801.26, 1),
    @(2021, "EricNC",   12,   5, 1819.48,             1815.72, 0),
    @(2021, "ChrisNC",  5,    9, 1861.72,             1850.76, 1),
    @(2021, "Erik",     1.5,  11, 2007.14,            1773.96, 1),
    @(2021, "Marcus",   1.5,  11, 2077.1999999999998, 1815.92, 1),
    @(2021, "Alex",     11,   1, 1495.16,             2031.4,  0)
)

$startRow = 126
$endRow = $startRow + $newRows.Count - 1

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $data = $newRows[$i]

    for ($col = 1; $col -le 7; $col++) {
        $cell = $ws.Cells.Item($row, $col)
        $cell.Value = $data[$col - 1]
        $cell.HorizontalAlignment = -4108
    }
}

# Extend column J's helper formula down through the new rows. Re-asserting
# the existing master formula across the full J67:J137 block keeps the
# untouched rows (67-125) exactly as they were and fills in the new ones.
$ws.Range("J67:J" + $endRow).Formula = "=IF(C67=1,1,0)"
$ws.Range("J" + $startRow + ":J" + $endRow).HorizontalAlignment = -4108

# Reflect the user's final on-screen position: scrolled down to the new
# rows with F138 (just past the last entry) selected.
$ws.Range("F" + ($endRow + 1)).Select()
